$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '25.260.52'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -2.86%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.550.87'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -4.94%  '

$ws.Range("E4").Value = '  -0.07%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '206.61'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -3.66%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.477'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -5.39%  '

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.0609'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.78%  '

$ws.Range("E9").Value = '  -3.27%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '17.75'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -4.27%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0778'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.42%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.766.28'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -4.92%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '1.544.16'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -5.34%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '3.98'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -4.92%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '0.504'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -4.85%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '25.252.25'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.90%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.0₃0707'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -4.66%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '58.59'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -4.82%  '

$ws.Range("E19").Value = '  -0.09%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '185.51'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -4.27%  '

$ws.Range("E21").Value = '  -3.80%  '

$ws.Range("E22").Value = '  -3.29%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '5.82'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -4.38%  '

$ws.Range("E24").Value = '  -4.08%  '

$ws.Range("E25").Value = '  -0.07%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '139.04'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -3.54%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '14.82'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -3.08%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '6.38'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -5.57%  '

$ws.Range("E30").Value = '  -6.64%  '

$ws.Range("E31").Value = '  -4.47%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '3.02'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -3.74%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '2.96'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -5.11%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '1.45'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -3.17%  '

$ws.Range("E35").Value = '  -4.01%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '1.084.48'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -3.62%  '

$ws.Range("E37").Value = '  -0.44%  '

$ws.Range("E38").Value = '  -2.79%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.492'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -5.87%  '

$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.24'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -7.82%  '

$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.759'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -10.98%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.798'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +3.88%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '92.68'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -5.66%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.681.42'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -4.87%  '

$ws.Range("E46").Value = '  -5.61%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.45'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -1.97%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '52.24'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -4.30%  '

$ws.Range("E49").Value = '  -5.70%  '

$ws.Range("E50").Value = '  -0.15%  '

$ws.Range("E51").Value = '  -2.20%  '
